# Updates targets-template.xlsx: refreshes the industry/region/
# document/form-type/visibility enumerations used by the data
# validation dropdowns (ValidationData sheet), the dropdown source
# ranges themselves, and the matching cell comments on row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("targets template")
$vd = $wb.Worksheets.Item("ValidationData")

$colG = @("Software", "Fintech", "Healthcare", "Medtech", "AI", "Computing", "Deep tech", "Climate", "Consumer", "E-commerce", "Marketplace", "Gaming", "Web3", "Developer tools", "Cybersecurity", "Logistics", "Adtech", "Proptech", "Agriculture", "Automotive", "Biotechnology", "Construction", "Education", "Energy", "Entertainment", "Environment", "Fashion", "Real estate", "Food", "IoT", "Government", "Hospitality", "HR", "Insurance", "Security", "Social", "Aerospace", "AR/VR", "Mining", "Advanced Materials", "Biofuels", "Hardware", "Nanotechnology", "Legal", "Manufacturing", "Media", "Pharmaceuticals", "Retail", "Telecommunications", "Transportation", "Agnostic", "Other")
$colH = @("Global", "North America", "South America", "LATAM", "Europe", "Middle East", "Africa", "Asia", "East Asia", "South East Asia", "South Asia", "Oceania", "EMEA", "Emerging Markets", "India", "China", "Japan", "Korea", "Australia", "United States", "Canada", "UK", "France", "Nigeria", "Kenya", "Egypt", "Senegal", "South Africa", "Netherlands", "Sweden", "Other")
$colK = @("pitch_deck", "video", "financials", "business_plan")
$colN = @("contact", "airtable", "typeform", "google", "generic")
$colP = @("FREE", "PRO", "MAX", "ENTERPRISE")

$oldCounts = @{ "G" = 45; "H" = 17; "K" = 5; "N" = 4; "P" = 3 }
$newCols = @{ "G" = $colG; "H" = $colH; "K" = $colK; "N" = $colN; "P" = $colP }

foreach ($col in @("G", "H", "K", "N", "P")) {
    $values = $newCols[$col]
    $oldCount = $oldCounts[$col]

    for ($i = 0; $i -lt $values.Length; $i++) {
        $addr = "$col" + ($i + 1)
        $vd.Range($addr).Value = $values[$i]
    }

    if ($oldCount -gt $values.Length) {
        for ($i = $values.Length + 1; $i -le $oldCount; $i++) {
            $addr = "$col" + $i
            $vd.Range($addr).ClearContents()
        }
    }

    $newCount = $values.Length
    $formula = "=ValidationData!`$" + $col + "`$1:`$" + $col + "`$" + $newCount
    $ws.Range("$col" + "2:" + "$col" + "9").Validation.Formula1 = $formula
    $ws.Range("$col" + "10:" + "$col" + "1000").Validation.Formula1 = $formula
}

# Refresh the "Possible values" comments so they mirror the new lists
$commentG = @"
Possible values:

- Software
- Fintech
- Healthcare
- Medtech
- AI
- Computing
- Deep tech
- Climate
- Consumer
- E-commerce
- Marketplace
- Gaming
- Web3
- Developer tools
- Cybersecurity
- Logistics
- Adtech
- Proptech
- Agriculture
- Automotive
- Biotechnology
- Construction
- Education
- Energy
- Entertainment
- Environment
- Fashion
- Real estate
- Food
- IoT
- Government
- Hospitality
- HR
- Insurance
- Security
- Social
- Aerospace
- AR/VR
- Mining
- Advanced Materials
- Biofuels
- Hardware
- Nanotechnology
- Legal
- Manufacturing
- Media
- Pharmaceuticals
- Retail
- Telecommunications
- Transportation
- Agnostic
- Other
"@
$null = $ws.Range("G1").Comment.Text($commentG)

$commentH = @"
Possible values:

- Global
- North America
- South America
- LATAM
- Europe
- Middle East
- Africa
- Asia
- East Asia
- South East Asia
- South Asia
- Oceania
- EMEA
- Emerging Markets
- India
- China
- Japan
- Korea
- Australia
- United States
- Canada
- UK
- France
- Nigeria
- Kenya
- Egypt
- Senegal
- South Africa
- Netherlands
- Sweden
- Other
"@
$null = $ws.Range("H1").Comment.Text($commentH)

$commentK = @"
Possible values:

- pitch_deck
- video
- financials
- business_plan
"@
$null = $ws.Range("K1").Comment.Text($commentK)

$commentN = @"
Possible values:

- contact
- airtable
- typeform
- google
- generic
"@
$null = $ws.Range("N1").Comment.Text($commentN)

$commentP = @"
Possible values:

- FREE
- PRO
- MAX
- ENTERPRISE
"@
$null = $ws.Range("P1").Comment.Text($commentP)

